$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and 1h-volume-change figures to the latest scrape.
# Cells are stored as text (not numbers), so force Text format before writing
# the new value to avoid Excel reinterpreting "302.03" / "-0.60%" as numeric/percent.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.024"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.62%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07827"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.35%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.187"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.98%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.08%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.033"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.74%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9136"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.06%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09717"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1890"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.69%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08580"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.47%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03525"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.17%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09976"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.75%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001481"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.19%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005642"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.89%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.468"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.076"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.41%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.70%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1301"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.759"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "10.57%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2205"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.04%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04634"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.29%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.93%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004798"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.31%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004751"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "28.39%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.52%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04729"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.63%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008069"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1393"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007666"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.47%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01043"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "13.35%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006050"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.61%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.613"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "142.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
